# The post "「お母さんがこの方法で僕を運ぶ時、暴れる余地はない」" (row 209) was
# removed from the source data. Delete its entire row so every subsequent
# row shifts up by one, matching the renumbered rows 209-362 in the target
# workbook (old dimension A1:C363 -> new A1:C362).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(209).Delete()
